$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "fastq_reformat/" subdirectory segment from each fastqFileName
# path stored in column F (rows 2-14).
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $old = $cell.Value()
    $new = $old -replace "fastq_reformat/", ""
    $cell.Value = $new
}

# Move the active selection to F15 (matches the post-edit selection in the file).
$ws.Range("F15").Select()
